# Update DBHAPI_Variabler Excel with new filter (Variabelliste) and
# group by (Group_by) values for tables 704, 705, 706, 707, 123, 124.
# This adds "Studiumkode" into the variable/group-by lists for these rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tables 704, 705, 706, 707 (Gjennomføring og frafall ...):
#   Variabelliste (E): "Institusjonskode,Avdkode,Årstall"
#                    -> "Institusjonskode,Avdkode, Studiumkode,Årstall"
#   Group_by (F):      "Institusjonskode,Avdkode,Årstall, Semester,Kjønnkode"
#                    -> "Institusjonskode,Avdkode, Studiumkode,Årstall, Semester,Kjønnkode"
$rows704to707 = 4, 5, 6, 7
foreach ($r in $rows704to707) {
    $ws.Cells.Item($r, 5).Value = "Institusjonskode,Avdkode, Studiumkode,Årstall"
    $ws.Cells.Item($r, 6).Value = "Institusjonskode,Avdkode, Studiumkode,Årstall, Semester,Kjønnkode"
}

# Table 123 (Registrerte studenter), row 36:
#   Variabelliste (E): "Institusjonskode,Avdelingskode,Årstall"
#                    -> "Institusjonskode,Avdelingskode, Studiumkode,Årstall"
#   Group_by (F):      "Institusjonskode,Avdelingskode,Årstall,Semester,Studentkategori"
#                    -> "Institusjonskode,Avdelingskode,Studiumkode,Årstall,Semester,Studentkategori"
$ws.Cells.Item(36, 5).Value = "Institusjonskode,Avdelingskode, Studiumkode,Årstall"
$ws.Cells.Item(36, 6).Value = "Institusjonskode,Avdelingskode,Studiumkode,Årstall,Semester,Studentkategori"

# Table 124 (Registrerte studenter fordelt på studieretning og campus), row 37:
#   Variabelliste (E): "Institusjonskode,Avdelingskode,Årstall,Stedkode campus"
#                    -> "Institusjonskode,Avdelingskode,Studiumkode,Årstall,Stedkode campus"
#   Group_by (F):      "Institusjonskode,Avdelingskode,Årstall,Semester,Studentkategori,Studieprogramkode,Studieretningkode,Stedkode campus"
#                    -> "Institusjonskode,Avdelingskode,Studiumkode,Årstall,Semester,Studentkategori,Studieprogramkode,Studieretningkode,Stedkode campus"
$ws.Cells.Item(37, 5).Value = "Institusjonskode,Avdelingskode,Studiumkode,Årstall,Stedkode campus"
$ws.Cells.Item(37, 6).Value = "Institusjonskode,Avdelingskode,Studiumkode,Årstall,Semester,Studentkategori,Studieprogramkode,Studieretningkode,Stedkode campus"

# Reflect the scrolled/selected view state left by the edit session.
$ws.Application.ActiveWindow.ScrollRow = 23
$ws.Range("E34").Select()
